$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at the top (shifts everything down by 1)
$ws.Rows.Item(1).Insert()

# 2. Clear column A's numeric index values (now at A3:A8, header A2 was already blank)
$ws.Range("A2:A8").ClearContents()

# 3. Set row heights for rows 1-8 to 27 points
$ws.Range("A1:A8").EntireRow.RowHeight = 27

# 4. Column widths (best effort; A and B are the meaningfully sized columns)
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 15

Write-Output "step done"
